# "Generate Report for Handback"
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for the two content files on the zh-cn and de-de sheets,
# flips the Status everywhere from "In Translation" to the handed-back message,
# and widens a few columns so the new/longer values are readable.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR long for RGB(0x64,0x95,0xED) == OOXML color FF6495ED

$targetUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5e7c7fda64633e12b15350da7daab5da9abbdf5/e2e/8ef2c0d9-a60e-4867-a4c4-345f00a1a6b0.md"
$targetUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5e7c7fda64633e12b15350da7daab5da9abbdf5/e2e/e2354d22-27a0-4251-a87c-91478321edea.md"
$display1 = "8ef2c0d9-a60e-4867-a4c4-345f00a1a6b0.md"
$display2 = "e2354d22-27a0-4251-a87c-91478321edea.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn / de-de) for both rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Row 2 - 8ef2c0d9... file
$wsZh.Range("I2").Value = $display1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetUrl1, "", "", $display1)
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("J2").Value = "8ef2c0d9-a60e-4867-a4c4-345f00a1a6b0.9c66d13bd03b5d878f1eb32052ead173117b16fe.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-24 10:13:27"

# Row 3 - e2354d22... file
$wsZh.Range("I3").Value = $display2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetUrl2, "", "", $display2)
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = $hyperlinkColor
$wsZh.Range("J3").Value = "e2354d22-27a0-4251-a87c-91478321edea.bc0418ec981e6cd0a6e0a0ff94d4c378a1b271d1.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-24 10:13:27"

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Row 2 - 8ef2c0d9... file
$wsDe.Range("I2").Value = $display1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetUrl1, "", "", $display1)
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("J2").Value = "8ef2c0d9-a60e-4867-a4c4-345f00a1a6b0.9c66d13bd03b5d878f1eb32052ead173117b16fe.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-24 10:13:43"

# Row 3 - e2354d22... file
$wsDe.Range("I3").Value = $display2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetUrl2, "", "", $display2)
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = $hyperlinkColor
$wsDe.Range("J3").Value = "e2354d22-27a0-4251-a87c-91478321edea.bc0418ec981e6cd0a6e0a0ff94d4c378a1b271d1.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-24 10:13:43"

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
